$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Replace the buggy "StatQuery" (column C) for the Cases/Samples/Files rows
# with the corrected program/study/case/sample/file counting query. All three
# rows previously shared the same broken Cypher text (shared-string 10); the
# fix swaps it for the new query the commit introduces, fixing ICDC breed
# testcases per the commit message.
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Chihuahua']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# --- View-state touch-ups captured in the sheetView: the author scrolled the
# sheet down a couple of rows and zoomed in from 55% to 85%, landing the
# selection on B4.
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws.Range("B4").Select() | Out-Null
